$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("total" -> shift to C, etc.)
# This shifts B:H to C:I, preserving values/styles/number formats.
$ws.Columns.Item(2).Insert()

# New header for inserted column B
$ws.Range("B1").Value = "Time (min)"

# Update data rows 2-19 with final values (values changed, not just shifted)
# Row 2: llama3.2_llama3.2
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 17.00299503588317
$ws.Range("E2").Value = 60.60333333333333
$ws.Range("F2").Value = 7.486666666666667
$ws.Range("G2").Value = 0.7325412531693777
$ws.Range("H2").Value = 0.1764804465224213
$ws.Range("I2").Value = 0.3838666666666667

# Row 3: llama3.1:8b_llama3.1:8b
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 11.99707562395166
$ws.Range("E3").Value = 57.96333333333333
$ws.Range("F3").Value = 7.43
$ws.Range("G3").Value = 0.7194825708866119
$ws.Range("H3").Value = 0.3042228413004278
$ws.Range("I3").Value = 0.4728166666666667

# Row 4: mistral:7b_mistral:7b
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 14.72542807699743
$ws.Range("E4").Value = 38.6275
$ws.Range("F4").Value = 9.7225
$ws.Range("G4").Value = 0.6809182316064835
$ws.Range("H4").Value = 0.2482662398409837
$ws.Range("I4").Value = 0.460025

# Row 5: gpt-3.5-turbo_gpt-3.5-turbo
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 14.54054299183933
$ws.Range("E5").Value = 56.838
$ws.Range("F5").Value = 9.344
$ws.Range("G5").Value = 0.7493282198905945
$ws.Range("H5").Value = 0.6377814874293684
$ws.Range("I5").Value = 0.53841

# Row 6: gpt-4o-mini_gpt-4o-mini
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 18.66998850620627
$ws.Range("E6").Value = 62.75
$ws.Range("F6").Value = 8.185
$ws.Range("G6").Value = 0.7023629918694496
$ws.Range("H6").Value = 0.2382443290059142
$ws.Range("I6").Value = 0.4912749999999999

# Row 7: gpt-4o_gpt-4o
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 12.58934115278937
$ws.Range("E7").Value = 63.475
$ws.Range("F7").Value = 7.148333333333333
$ws.Range("G7").Value = 0.796918253103892
$ws.Range("H7").Value = 0.165238055902091
$ws.Range("I7").Value = 0.448

# Row 8: llama3.2_llama3.2
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 34.7643808065357
$ws.Range("E8").Value = 55.45857142857143
$ws.Range("F8").Value = 9.705714285714285
$ws.Range("G8").Value = 0.8355420402118138
$ws.Range("H8").Value = 0.3524402083705272
$ws.Range("I8").Value = 0.5774

# Row 9: llama3.1:8b_llama3.1:8b
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 10.30424690887918
$ws.Range("E9").Value = 51.73666666666667
$ws.Range("F9").Value = 8.133333333333333
$ws.Range("G9").Value = 0.726477692524592
$ws.Range("H9").Value = 0.1034525033367635
$ws.Range("I9").Value = 0.44335

# Row 10: mistral:7b_mistral:7b
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 16.84469795342042
$ws.Range("E10").Value = 39.04166666666666
$ws.Range("F10").Value = 11.40833333333333
$ws.Range("G10").Value = 0.6543693641821543
$ws.Range("H10").Value = 0.6220974346825636
$ws.Range("I10").Value = 0.6211666666666668

# Row 11: gpt-3.5-turbo_gpt-3.5-turbo
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 15.15746019551164
$ws.Range("E11").Value = 54.8738888888889
$ws.Range("F11").Value = 10.63944444444445
$ws.Range("G11").Value = 0.7475307981173197
$ws.Range("H11").Value = 0.7030727832010455
$ws.Range("I11").Value = 0.6319111111111111

# Row 12: gpt-4o-mini_gpt-4o-mini
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 18.97582117255906
$ws.Range("E12").Value = 59.59
$ws.Range("F12").Value = 9.383999999999999
$ws.Range("G12").Value = 0.8014606833457947
$ws.Range("H12").Value = 0.4616497657644622
$ws.Range("I12").Value = 0.5099600000000001

# Row 13: gpt-4o_gpt-4o
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 16.68593439135956
$ws.Range("E13").Value = 58.26416666666666
$ws.Range("F13").Value = 10.1125
$ws.Range("G13").Value = 0.7464624593655268
$ws.Range("H13").Value = 0.2561081290134884
$ws.Range("I13").Value = 0.5094583333333333

# Row 14: llama3.2_llama3.2
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 18
$ws.Range("D14").Value = 9.354589395536147
$ws.Range("E14").Value = 47.06666666666668
$ws.Range("F14").Value = 11.845
$ws.Range("G14").Value = 0.8879779842164781
$ws.Range("H14").Value = 0.3910482174231577
$ws.Range("I14").Value = 0.6651222222222223

# Row 15: llama3.1:8b_llama3.1:8b
$ws.Range("B15").Value = 15
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 14.56586365409762
$ws.Range("E15").Value = 49.005
$ws.Range("F15").Value = 10.12
$ws.Range("G15").Value = 0.7706534206867218
$ws.Range("H15").Value = 0.1601834382184749
$ws.Range("I15").Value = 0.53461

# Row 16: mistral:7b_mistral:7b
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 13.95205316602343
$ws.Range("E16").Value = 31.89875
$ws.Range("F16").Value = 13.3475
$ws.Range("G16").Value = 0.6983311623334885
$ws.Range("H16").Value = 0.05448345582181986
$ws.Range("I16").Value = 0.56465

# Row 17: gpt-3.5-turbo_gpt-3.5-turbo
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = 28
$ws.Range("D17").Value = 15.60235699510913
$ws.Range("E17").Value = 53.46071428571427
$ws.Range("F17").Value = 11.2975
$ws.Range("G17").Value = 0.6930893010326794
$ws.Range("H17").Value = 0.2582317777198346
$ws.Range("I17").Value = 0.6400535714285713

# Row 18: gpt-4o-mini_gpt-4o-mini
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 17.14890865168301
$ws.Range("E18").Value = 66.72714285714285
$ws.Range("F18").Value = 8.34142857142857
$ws.Range("G18").Value = 0.80890120778765
$ws.Range("H18").Value = 0.3165135501517266
$ws.Range("I18").Value = 0.5204214285714286

# Row 19: gpt-4o_gpt-4o
$ws.Range("B19").Value = 15
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 18.45210486254826
$ws.Range("E19").Value = 62.67642857142857
$ws.Range("F19").Value = 9.41642857142857
$ws.Range("G19").Value = 0.7904890860830035
$ws.Range("H19").Value = 0.5723187685291828
$ws.Range("I19").Value = 0.5388785714285714

# Add new rows 20-24
# Row 20
$ws.Range("A20").Value = "gemma3:12b_gemma3:12b"
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 23.31019668279999
$ws.Range("E20").Value = 58.685
$ws.Range("F20").Value = 7.011666666666666
$ws.Range("G20").Value = 0.6479067802429199
$ws.Range("H20").Value = 0.370180423848563
$ws.Range("I20").Value = 0.4172166666666666

# Row 21
$ws.Range("A21").Value = "gemma3:12b_gemma3:12b"
$ws.Range("B21").Value = 10
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 17.180115432898
$ws.Range("E21").Value = 63.185
$ws.Range("F21").Value = 5.761666666666667
$ws.Range("G21").Value = 0.6554248730341593
$ws.Range("H21").Value = 0.1856545144847652
$ws.Range("I21").Value = 0.4137

# Row 22
$ws.Range("A22").Value = "gemma3:12b_gemma3:12b"
$ws.Range("B22").Value = 15
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 28.06364570366047
$ws.Range("E22").Value = 52.148
$ws.Range("F22").Value = 9.379000000000001
$ws.Range("G22").Value = 0.7058762311935425
$ws.Range("H22").Value = 0.3776658643586975
$ws.Range("I22").Value = 0.4523499999999999

# Row 23
$ws.Range("A23").Value = "gemini-2.5-pro-exp_gemini-2.5-pro-exp"
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 24.89709997424184
$ws.Range("E23").Value = 61.90666666666667
$ws.Range("F23").Value = 5.893333333333334
$ws.Range("G23").Value = 0.5419649879137675
$ws.Range("H23").Value = 0.2892044354747683
$ws.Range("I23").Value = 0.5014666666666666

# Row 24
$ws.Range("A24").Value = "gemini-2.5-pro-exp_gemini-2.5-pro-exp"
$ws.Range("B24").Value = 10
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 30.77090886883413
$ws.Range("E24").Value = 54.642
$ws.Range("F24").Value = 9.209
$ws.Range("G24").Value = 0.6478052794933319
$ws.Range("H24").Value = 0.4877557687088286
$ws.Range("I24").Value = 0.50371

